$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) - copy style from existing header cell (H1) so new headers match formatting
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for I and J columns (rows 2-15)
$data = @(
    @(6, 8),
    @(1, 3),
    @(6, 7),
    @(7, 9),
    @(2, 4),
    @(6, 6),
    @(5, 8),
    @(7, 8),
    @(8, 9),
    @(1, 4),
    @(1, 3),
    @(1, 3),
    @(1, 2),
    @(3, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
